# Applies the edits described in the commit:
#  - Sets cell D16 on Sheet1 to the text "Regressa a 1"
#  - Moves the active selection to D16 (matching the saved view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D16").Value = "Regressa a 1"

$ws.Activate()
$ws.Range("D16").Select()
